# Update with restock suggestion
# Applies the forecast_summary_B09JZGTYXJ_WITH_PO.xlsx edits:
#  - Sheet "Forecast Comparison": fill in Week_Start_Date (col B), refresh
#    Inventory Coverage (L), Stockout Risk (M), Reorder Urgency (N) and
#    Seasonality Index (P) per-row, drop the "Sales Volume Rank" column and
#    shift "Lifecycle Stage" (formerly R) left into Q with updated values.
#  - Sheet "Summary": refresh Max Forecast / Max Forecast Week / Min Forecast
#    Week metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

function Set-TextValue($range, $text) {
    # Force a literal text value (avoids Excel auto-converting things that
    # look like dates/numbers into date serials / numbers), then reset the
    # cell style back to Normal so no stray NumberFormat/style sticks around.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Week_Start_Date (column B) ---------------------------------------
Set-TextValue $ws.Range("B2")  "2025-02-02"
Set-TextValue $ws.Range("B3")  "2025-02-09"
Set-TextValue $ws.Range("B4")  "2025-02-16"
Set-TextValue $ws.Range("B5")  "2025-02-23"
Set-TextValue $ws.Range("B6")  "2025-03-02"
Set-TextValue $ws.Range("B7")  "2025-03-09"
Set-TextValue $ws.Range("B8")  "2025-03-16"
Set-TextValue $ws.Range("B9")  "2025-03-23"
Set-TextValue $ws.Range("B10") "2025-03-30"
Set-TextValue $ws.Range("B11") "2025-04-06"
Set-TextValue $ws.Range("B12") "2025-04-13"
Set-TextValue $ws.Range("B13") "2025-04-20"
Set-TextValue $ws.Range("B14") "2025-04-27"
Set-TextValue $ws.Range("B15") "2025-05-04"
Set-TextValue $ws.Range("B16") "2025-05-11"
Set-TextValue $ws.Range("B17") "2025-05-18"

# --- Inventory Coverage (L), Stockout Risk (M), Reorder Urgency (N), ---
# --- Seasonality Index (P) ---------------------------------------------
$ws.Range("L2").Value = 2.5
$ws.Range("M2").Value = "Low"
$ws.Range("N2").Value = "Normal"
$ws.Range("P2").Value = 1.18

$ws.Range("L3").Value = 1.5
$ws.Range("M3").Value = "Low"
$ws.Range("N3").Value = "Normal"
$ws.Range("P3").Value = 0.8

$ws.Range("L4").Value = 0.77
$ws.Range("M4").Value = "Low"
$ws.Range("P4").Value = 0.99

$ws.Range("L5").Value = 0
$ws.Range("P5").Value = 1.03

$ws.Range("L6").Value = 0
$ws.Range("P6").Value = 1.04

$ws.Range("L7").Value = 0
$ws.Range("P7").Value = 0.93

$ws.Range("L8").Value = 0
$ws.Range("P8").Value = 0.84

$ws.Range("L9").Value = 0
$ws.Range("P9").Value = 0.9

$ws.Range("L10").Value = 0
$ws.Range("P10").Value = 1.13

$ws.Range("L11").Value = 0
$ws.Range("P11").Value = 0.9

$ws.Range("L12").Value = 0
$ws.Range("P12").Value = 1.09

$ws.Range("L13").Value = 0
$ws.Range("P13").Value = 1.09

$ws.Range("L14").Value = 0
$ws.Range("P14").Value = 1.06

$ws.Range("L15").Value = 0
$ws.Range("P15").Value = 1.13

$ws.Range("L16").Value = 0
$ws.Range("P16").Value = 1.12

$ws.Range("L17").Value = 0
$ws.Range("P17").Value = 0.9399999999999999

# --- Drop "Sales Volume Rank" (Q) and shift "Lifecycle Stage" (R) left -
# Every row's Lifecycle Stage becomes "Mature" (was "Growth"); move that
# column from R into Q, then delete the now-empty R column entirely.
$ws.Range("Q1").Value = "Lifecycle Stage"
for ($r = 2; $r -le 17; $r++) {
    $ws.Range("Q$r").Value = "Mature"
}
$ws.Columns.Item(18).Delete()

# --- Summary sheet -------------------------------------------------------
# Column B on this sheet stores every metric as text (even numeric-looking
# ones), so force "1" to stay text instead of becoming a number.
$summary = $wb.Worksheets.Item("Summary")
Set-TextValue $summary.Range("B12") "1"
$summary.Range("B13").Value = "N/A"
$summary.Range("B15").Value = "N/A"
